# GroupTagTests_TotalLabel.xlsx template fix
#
# ClosedXML.Report was bumped to 0.100.3. That upgrade made it impossible to
# stash a "display value" inside a cell formula (e.g. &="Total: "<<sum>>) and
# later read back the original templated text - the formula is now actually
# evaluated instead. The template's "Formulas in group row" sample cell
# (G6) relied on that old (incorrect) behaviour, so it's switched to a plain
# templated text tag that does not use a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G6 used to hold the literal text '&="Total: "<<sum>>' (a fake formula used
# only to stash a display string). Replace it with a plain <<sum>> tag.
$ws.Range("G6").Value = "<<sum>>"
